$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.980.68'
$ws.Range("E2").Value = '  +2.25%  '
$ws.Range("D3").Value = '3.190.33'
$ws.Range("E3").Value = '  +1.16%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Formula = "'535.93"
$ws.Range("E5").Value = '  +0.83%  '
$ws.Range("D6").Formula = "'145.02"
$ws.Range("E6").Value = '  +3.84%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -1.34%  '
$ws.Range("D9").Formula = "'7.32"
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("E10").Value = '  +1.49%  '
$ws.Range("D11").Formula = "'0.431"
$ws.Range("E11").Value = '  -0.97%  '
$ws.Range("D12").Value = '3.743.18'
$ws.Range("E12").Value = '  +1.23%  '
$ws.Range("E13").Value = '  -2.21%  '
$ws.Range("D14").Formula = "'25.79"
$ws.Range("E14").Value = '  -0.17%  '
$ws.Range("E15").Value = '  +0.25%  '
$ws.Range("D16").Value = '60.009.95'
$ws.Range("E16").Value = '  +2.22%  '
$ws.Range("D17").Value = '3.196.09'
$ws.Range("E17").Value = '  +1.30%  '
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("D19").Formula = "'13.25"
$ws.Range("E19").Value = '  +1.99%  '
$ws.Range("D20").Formula = "'8.19"
$ws.Range("E20").Value = '  +0.60%  '
$ws.Range("D21").Formula = "'369.30"
$ws.Range("E21").Value = '  -0.65%  '
$ws.Range("D22").Formula = "'0.998"
$ws.Range("E22").Value = '  -0.25%  '
$ws.Range("E23").Value = '  -0.91%  '
$ws.Range("D24").Formula = "'69.35"
$ws.Range("E24").Value = '  -0.44%  '
$ws.Range("E25").Value = '  +1.01%  '
$ws.Range("E26").Value = '  +3.66%  '
$ws.Range("E27").Value = '  -0.48%  '
$ws.Range("D28").Value = '0.0₃0872'
$ws.Range("E28").Value = '  +0.67%  '
$ws.Range("D29").Formula = "'22.47"
$ws.Range("E29").Value = '  +1.43%  '
$ws.Range("E30").Value = '  +0.53%  '
$ws.Range("D31").Formula = "'6.10"
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("D32").Formula = "'5.29"
$ws.Range("E32").Value = '  +2.46%  '
$ws.Range("E33").Value = '  +2.55%  '
$ws.Range("D34").Formula = "'6.57"
$ws.Range("E34").Value = '  +4.03%  '
$ws.Range("D35").Formula = "'156.55"
$ws.Range("E35").Value = '  -1.43%  '
$ws.Range("E36").Value = '  +1.81%  '
$ws.Range("D37").Formula = "'26.60"
$ws.Range("E37").Value = '  +5.90%  '
$ws.Range("D38").Value = '2.796.00'
$ws.Range("E38").Value = '  +5.94%  '
$ws.Range("D39").Formula = "'0.0706"
$ws.Range("E39").Value = '  +2.89%  '
$ws.Range("D40").Formula = "'0.0308"
$ws.Range("E40").Value = '  +7.57%  '
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("D42").Formula = "'4.25"
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("D43").Formula = "'39.63"
$ws.Range("E43").Value = '  +1.48%  '
$ws.Range("E44").Value = '  +1.64%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").Formula = "'0.104"
$ws.Range("E45").Value = '  +0.54%  '
$ws.Range("B46").Value = 'RenzoRestakedETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D46").Value = '3.232.92'
$ws.Range("E46").Value = '  +1.19%  '
$ws.Range("D47").Formula = "'0.983"
$ws.Range("E47").Value = '  +0.31%  '
$ws.Range("D48").Formula = "'6.14"
$ws.Range("E48").Value = '  -0.95%  '
$ws.Range("D49").Formula = "'20.66"
$ws.Range("E49").Value = '  +2.08%  '
$ws.Range("D50").Formula = "'0.792"
$ws.Range("E50").Value = '  +4.76%  '
